$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to keep their text formatting (e.g. trailing
# zeros like "1.00" or "0.510") instead of Excel auto-converting them to numbers.
$priceCells = @(
    "D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D18",
    "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D30", "D32", "D33", "D36",
    "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.335.31"
$ws.Range("E2").Value = "  +5.12%  "

# Row 3
$ws.Range("D3").Value = "3.246.19"
$ws.Range("E3").Value = "  +2.72%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "577.12"
$ws.Range("E5").Value = "  +2.63%  "

# Row 6
$ws.Range("D6").Value = "179.14"
$ws.Range("E6").Value = "  +6.47%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -0.35%  "

# Row 9
$ws.Range("D9").Value = "3.240.42"
$ws.Range("E9").Value = "  +2.56%  "

# Row 10
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +4.69%  "

# Row 11
$ws.Range("D11").Value = "6.74"
$ws.Range("E11").Value = "  +2.07%  "

# Row 12
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +4.82%  "

# Row 13
$ws.Range("D13").Value = "3.803.22"
$ws.Range("E13").Value = "  +2.59%  "

# Row 14
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").Value = "27.92"
$ws.Range("E15").Value = "  +2.40%  "

# Row 16
$ws.Range("D16").Value = "67.236.03"
$ws.Range("E16").Value = "  +4.98%  "

# Row 17
$ws.Range("E17").Value = "  +3.31%  "

# Row 18
$ws.Range("D18").Value = "3.246.32"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  +2.47%  "

# Row 20
$ws.Range("D20").Value = "13.36"
$ws.Range("E20").Value = "  +3.43%  "

# Row 21
$ws.Range("D21").Value = "375.86"
$ws.Range("E21").Value = "  +6.93%  "

# Row 22
$ws.Range("D22").Value = "7.60"
$ws.Range("E22").Value = "  +6.39%  "

# Row 23
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("D24").Value = "71.25"
$ws.Range("E24").Value = "  +3.85%  "

# Row 25
$ws.Range("D25").Value = "0.510"
$ws.Range("E25").Value = "  +1.84%  "

# Row 26
$ws.Range("D26").Value = "3.387.56"
$ws.Range("E26").Value = "  +2.50%  "

# Row 27
$ws.Range("E27").Value = "  -0.70%  "

# Row 28
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  +5.03%  "

# Row 29
$ws.Range("E29").Value = "  +2.02%  "

# Row 30
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
$ws.Range("E31").Value = "  +4.71%  "

# Row 32
$ws.Range("D32").Value = "5.63"
$ws.Range("E32").Value = "  +2.22%  "

# Row 33
$ws.Range("D33").Value = "22.52"
$ws.Range("E33").Value = "  +2.83%  "

# Row 34
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("E35").Value = "  +6.29%  "

# Row 36
$ws.Range("D36").Value = "6.82"
$ws.Range("E36").Value = "  +3.07%  "

# Row 37
$ws.Range("D37").Value = "163.19"
$ws.Range("E37").Value = "  +6.30%  "

# Row 38
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +4.16%  "

# Row 39
$ws.Range("D39").Value = "0.856"
$ws.Range("E39").Value = "  +5.34%  "

# Row 40
$ws.Range("E40").Value = "  +9.96%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "6.83"
$ws.Range("E41").Value = "  +14.38%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "26.87"
$ws.Range("E42").Value = "  +4.73%  "

# Row 43
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  +4.21%  "

# Row 44
$ws.Range("D44").Value = "2.761.86"
$ws.Range("E44").Value = "  +6.37%  "

# Row 45
$ws.Range("D45").Value = "4.40"
$ws.Range("E45").Value = "  +5.60%  "

# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "352.53"
$ws.Range("E46").Value = "  +11.26%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "25.74"
$ws.Range("E47").Value = "  +9.02%  "

# Row 48
$ws.Range("D48").Value = "40.45"
$ws.Range("E48").Value = "  +2.81%  "

# Row 49
$ws.Range("D49").Value = "0.0673"
$ws.Range("E49").Value = "  +3.50%  "

# Row 50
$ws.Range("E50").Value = "  +4.48%  "

# Row 51
$ws.Range("E51").Value = "  +1.77%  "
